# Apply the tracked-progress updates to Sheet1: new SUM() formulas for the
# weekly "written" totals (B3/B4/B5), an updated daily page count (I15), and
# move the active-cell selection to I16 to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Preface week already totalled manually; now driven by the daily log.
$ws.Range("B3").Formula = "=SUM(I1:I4)"

# Introduction week: was a hard-coded 14, now summed from the daily log
# (I5:I16), which also reflects today's update to I15 below.
$ws.Range("B4").Formula = "=SUM(I5:I16)"

# Theory week: likewise now driven by the daily log.
$ws.Range("B5").Formula = "=SUM(I17:I21)"

# Today's page count.
$ws.Range("I15").Value = 6

# Leave the selection where the author last clicked.
$ws.Range("I16").Select() | Out-Null
